$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same for every data row, rows 2-26), columns B..I
$newValues = @{
    2 = 0.85251841622735       # B: r2
    3 = 0.8433531604744468     # C: r2_test
    4 = 0.950923143859281      # D: r2_val
    5 = 0.8803231797139723     # E: r2_vt
    6 = 0.1632186621427536     # F: mse
    7 = 0.3310821652412415     # G: mse_test
    8 = 0.05881708115339279    # H: mse_val
    9 = 0.2029581218957901     # I: mse_vt
}

for ($r = 2; $r -le 26; $r++) {
    foreach ($c in $newValues.Keys) {
        $ws.Cells.Item($r, $c).Value = $newValues[$c]
    }
}
